# Auto-generated Excel COM-interop script
# Applies the numeric corrections described in the commit diff
# for the 'Tonberry_Profits' workbook (per-sheet Leve profit tables).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1029.7307
$ws.Range("I15").Value = 1029.7307
$ws.Range("K15").Value = 3089.1921
$ws.Range("M15").Value = -2920.1921

# Row 62
$ws.Range("H62").Value = 6313.1665
$ws.Range("I62").Value = 1979.8
$ws.Range("J62").Value = 27980
$ws.Range("K62").Value = 1979.8
$ws.Range("L62").Value = 27980
$ws.Range("M62").Value = -1355.8
$ws.Range("N62").Value = -29228

# Row 65
$ws.Range("H65").Value = 6313.1665
$ws.Range("I65").Value = 1979.8
$ws.Range("J65").Value = 27980
$ws.Range("K65").Value = 9899
$ws.Range("L65").Value = 139900
$ws.Range("M65").Value = -6779
$ws.Range("N65").Value = -146140

# Row 80
$ws.Range("H80").Value = 2818.5
$ws.Range("I80").Value = 3169.2856
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 9507.856800000001
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -8509.856800000001
$ws.Range("N80").Value = -7996

# Row 83
$ws.Range("H83").Value = 2818.5
$ws.Range("I83").Value = 3169.2856
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 28523.5704
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -23531.5704
$ws.Range("N83").Value = -27984

# Row 94
$ws.Range("H94").Value = 1791.875
$ws.Range("I94").Value = 1791.875
$ws.Range("K94").Value = 1791.875
$ws.Range("M94").Value = -1340.875

# Row 98
$ws.Range("H98").Value = 741.37836
$ws.Range("I98").Value = 488.13333
$ws.Range("J98").Value = 1826.7142
$ws.Range("K98").Value = 488.13333
$ws.Range("L98").Value = 1826.7142
$ws.Range("M98").Value = 1009.86667
$ws.Range("N98").Value = -4822.7142

# Row 106
$ws.Range("H106").Value = 2497.5
$ws.Range("I106").Value = 2497.5
$ws.Range("K106").Value = 2497.5
$ws.Range("M106").Value = -1866.5

# Row 122
$ws.Range("H122").Value = 741.37836
$ws.Range("I122").Value = 488.13333
$ws.Range("J122").Value = 1826.7142
$ws.Range("K122").Value = 1464.39999
$ws.Range("L122").Value = 5480.142599999999
$ws.Range("M122").Value = 985.6000100000001
$ws.Range("N122").Value = -10380.1426

# Row 132
$ws.Range("H132").Value = 1594.421
$ws.Range("I132").Value = 1516.3334
$ws.Range("K132").Value = 4549.0002
$ws.Range("M132").Value = -2019.0002

# Row 137
$ws.Range("H137").Value = 2666.4443
$ws.Range("I137").Value = 2114.5
$ws.Range("J137").Value = 2942.4167
$ws.Range("K137").Value = 6343.5
$ws.Range("L137").Value = 8827.250100000001
$ws.Range("M137").Value = -3793.5
$ws.Range("N137").Value = -13927.2501

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 1005.6429
$ws.Range("I74").Value = 527.0571
$ws.Range("K74").Value = 527.0571
$ws.Range("M74").Value = 346.9429

# Row 77
$ws.Range("H77").Value = 1005.6429
$ws.Range("I77").Value = 527.0571
$ws.Range("K77").Value = 2635.2855
$ws.Range("M77").Value = 1732.7145

$ws = $wb.Worksheets.Item("BSM")
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 134
$ws.Range("H134").Value = 5862.8213
$ws.Range("I134").Value = 6671.7393
$ws.Range("K134").Value = 20015.2179
$ws.Range("M134").Value = -17480.2179

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 86440.42999999999
$ws.Range("I16").Value = 86440.42999999999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 86440.42999999999
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -86153.42999999999
$ws.Range("N16").ClearContents()

# Row 31
$ws.Range("H31").Value = 3906.2727
$ws.Range("I31").Value = 8493
$ws.Range("J31").Value = 2887
$ws.Range("K31").Value = 8493
$ws.Range("L31").Value = 2887
$ws.Range("M31").Value = -8198
$ws.Range("N31").Value = -3477

# Row 34
$ws.Range("H34").Value = 3906.2727
$ws.Range("I34").Value = 8493
$ws.Range("J34").Value = 2887
$ws.Range("K34").Value = 8493
$ws.Range("L34").Value = 2887
$ws.Range("M34").Value = -8291
$ws.Range("N34").Value = -3291

# Row 62
$ws.Range("H62").Value = 2600.8572
$ws.Range("I62").Value = 2617.8333
$ws.Range("J62").Value = 2499
$ws.Range("K62").Value = 2617.8333
$ws.Range("L62").Value = 2499
$ws.Range("M62").Value = -1993.8333
$ws.Range("N62").Value = -3747

# Row 65
$ws.Range("H65").Value = 2600.8572
$ws.Range("I65").Value = 2617.8333
$ws.Range("J65").Value = 2499
$ws.Range("K65").Value = 13089.1665
$ws.Range("L65").Value = 12495
$ws.Range("M65").Value = -9969.166499999999
$ws.Range("N65").Value = -18735

# Row 86
$ws.Range("H86").Value = 3023.7144
$ws.Range("I86").Value = 2433.2
$ws.Range("K86").Value = 2433.2
$ws.Range("M86").Value = -1310.2

# Row 89
$ws.Range("H89").Value = 3023.7144
$ws.Range("I89").Value = 2433.2
$ws.Range("K89").Value = 12166
$ws.Range("M89").Value = -6550

# Row 92
$ws.Range("H92").Value = 38999
$ws.Range("J92").Value = 38999
$ws.Range("L92").Value = 38999
$ws.Range("N92").Value = -43991

# Row 113
$ws.Range("H113").Value = 86440.42999999999
$ws.Range("I113").Value = 86440.42999999999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 86440.42999999999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -84270.42999999999
$ws.Range("N113").ClearContents()

# Row 134
$ws.Range("H134").Value = 4543.6665
$ws.Range("I134").Value = 4379.4
$ws.Range("J134").Value = 4749
$ws.Range("K134").Value = 13138.2
$ws.Range("L134").Value = 14247
$ws.Range("M134").Value = -10603.2
$ws.Range("N134").Value = -19317

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 151.125
$ws.Range("I40").Value = 173.75
$ws.Range("J40").Value = 128.5
$ws.Range("K40").Value = 695
$ws.Range("L40").Value = 514
$ws.Range("M40").Value = -626
$ws.Range("N40").Value = -652

# Row 131
$ws.Range("H131").Value = 16970.906
$ws.Range("J131").Value = 18637.154
$ws.Range("L131").Value = 55911.462
$ws.Range("N131").Value = -65991.462

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 4000
$ws.Range("I9").Value = 4000
$ws.Range("K9").Value = 4000
$ws.Range("M9").Value = -3830

# Row 97
$ws.Range("H97").Value = 1856.7273
$ws.Range("I97").Value = 1880.4445
$ws.Range("J97").Value = 1750
$ws.Range("K97").Value = 1880.4445
$ws.Range("L97").Value = 1750
$ws.Range("M97").Value = -1384.4445
$ws.Range("N97").Value = -2742

# Row 122
$ws.Range("H122").Value = 1696.7142
$ws.Range("I122").Value = 1576.4
$ws.Range("K122").Value = 4729.200000000001
$ws.Range("M122").Value = -2279.200000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

# Row 61
$ws.Range("H61").Value = 2686.3333
$ws.Range("I61").Value = 2217
$ws.Range("J61").Value = 3625
$ws.Range("K61").Value = 2217
$ws.Range("L61").Value = 3625
$ws.Range("M61").Value = -2015
$ws.Range("N61").Value = -4029

# Row 113
$ws.Range("H113").Value = 2686.3333
$ws.Range("I113").Value = 2217
$ws.Range("J113").Value = 3625
$ws.Range("K113").Value = 2217
$ws.Range("L113").Value = 3625
$ws.Range("M113").Value = -47
$ws.Range("N113").Value = -7965

# Row 132
$ws.Range("H132").Value = 2556.818
$ws.Range("I132").Value = 1458.1428
$ws.Range("K132").Value = 4374.428400000001
$ws.Range("M132").Value = -1844.428400000001

# Row 134
$ws.Range("H134").Value = 48426.5
$ws.Range("J134").Value = 48426.5
$ws.Range("L134").Value = 48426.5
$ws.Range("N134").Value = -58566.5

# Row 136
$ws.Range("H136").Value = 5463.1763
$ws.Range("I136").Value = 4090.6365
$ws.Range("K136").Value = 12271.9095
$ws.Range("M136").Value = -9721.9095

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 80005
$ws.Range("J7").Value = 80005
$ws.Range("L7").Value = 80005
$ws.Range("N7").Value = -80231

# Row 14
$ws.Range("H14").Value = 12000
$ws.Range("J14").Value = 12000
$ws.Range("L14").Value = 12000
$ws.Range("N14").Value = -12336

# Row 100
$ws.Range("H100").Value = 839.7
$ws.Range("I100").Value = 674.625
$ws.Range("K100").Value = 1349.25
$ws.Range("M100").Value = -808.25

# Row 133
$ws.Range("H133").Value = 61485
$ws.Range("J133").Value = 59998.332
$ws.Range("L133").Value = 59998.332
$ws.Range("N133").Value = -70118.33199999999

